# "Added 1.1.0 of term"
# Bump the Version and Date metadata values on the "Metadata" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# A3/B3 = "Version" / "1.0.0" -> "1.1.0"
$ws.Range("B3").Value = "1.1.0"

# A8/B8 = "Date" / "2023-06-07T11:52:14+02:00" -> "2023-07-10T23:08:03+02:00"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
